$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(1360.5, 1335.6, 1338.85, 1340.7, 58, 1345.3)
    3  = @(1175, 1156.25, 1164, 1164.6, 22, 1164.7)
    4  = @(48134.15, 47731, 47910, 47908.55, 19, 48059.95)
    5  = @(119.65, 116.9, 119.2, 119.2, 249, 117.15)
    6  = @(850.75, 827.7, 828, 829.3, 137, 831.6)
    7  = @(660.6, 651.7, 656.45, 657.15, 130, 657.75)
    8  = @(1131, 1118.8, 1129.65, 1129.05, 237, 1127.2)
    9  = @(1002.1, 978.7, 998.8, 998, 60, 991.25)
    10 = @(22375.45, 22211, 22280.5, 22287.05, 54, 22336)
    11 = @(2868.55, 2840.15, 2844.5, 2845.05, 48, 2850.6)
    12 = @(815.3, 807.65, 810.65, 810.55, 285, 812.55)
    13 = @(1082.5, 1064.2, 1068, 1067.3, 21, 1081.15)
    14 = @(973.55, 949.2, 949.6, 951.4, 243, 970.95)
    15 = @(168.45, 166, 166.5, 166.4, 652, 167.25)
    16 = @(3909.25, 3864, 3871.9, 3871.75, 14, 3905.9)
    17 = @(3323, 3274.05, 3285, 3284.95, 8, 3311.1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$ws.Range("I9").Select()
